$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") date value changes from 2023-10-13 (45212) to 2023-10-22 (45221)
# for every data row (rows 2 through 38).
for ($row = 2; $row -le 38; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
